# Auto-applied corrections to the ifrs company_list data rows (rows 2-9).
# Commit message: "error solve ifrs list" - fixes erroneous bulk financial
# figures and removes the now-unused "U" (and, for rows 7-9, Q/R/S/T) columns
# that were left over from a bad paste/merge.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 40138
$ws.Range("E2").Value = 2708
$ws.Range("F2").Value = 2708
$ws.Range("G2").Value = 2738
$ws.Range("H2").Value = 2058
$ws.Range("I2").Value = 2052
$ws.Range("J2").Value = 7
$ws.Range("K2").Value = 306133
$ws.Range("L2").Value = 264390
$ws.Range("M2").Value = 41743
$ws.Range("N2").Value = 41671
$ws.Range("O2").Value = 72
$ws.Range("P2").Value = 17039
$ws.Range("Q2").Value = 10918
$ws.Range("R2").Value = -30810
$ws.Range("S2").Value = 10372
$ws.Range("T2").Value = 53
$ws.Range("U2").ClearContents()
$ws.Range("V2").Value = 27439
$ws.Range("W2").Value = 6.75
$ws.Range("X2").Value = 5.13
$ws.Range("Y2").Value = 5.08
$ws.Range("Z2").Value = 0.73
$ws.Range("AA2").Value = 633.38
$ws.Range("AB2").Value = 151.84
$ws.Range("AC2").Value = 602
$ws.Range("AD2").Value = 15.96
$ws.Range("AE2").Value = 12440
$ws.Range("AF2").Value = 0.77
$ws.Range("AG2").Value = 250
$ws.Range("AH2").Value = 2.6
$ws.Range("AI2").Value = 40.96
$ws.Range("AJ2").Value = 326700891

# --- Row 3 ---
$ws.Range("D3").Value = 41925
$ws.Range("E3").Value = 1485
$ws.Range("F3").Value = 1485
$ws.Range("G3").Value = 2261
$ws.Range("H3").Value = 1746
$ws.Range("I3").Value = 1698
$ws.Range("J3").Value = 49
$ws.Range("K3").Value = 284160
$ws.Range("L3").Value = 249640
$ws.Range("M3").Value = 34520
$ws.Range("N3").Value = 34287
$ws.Range("O3").Value = 234
$ws.Range("P3").Value = 16981
$ws.Range("Q3").Value = -38435
$ws.Range("R3").Value = -36
$ws.Range("S3").Value = 41128
$ws.Range("T3").Value = 68
$ws.Range("U3").ClearContents()
$ws.Range("V3").Value = 29465
$ws.Range("W3").Value = 3.54
$ws.Range("X3").Value = 4.17
$ws.Range("Y3").Value = 4.47
$ws.Range("Z3").Value = 0.57
$ws.Range("AA3").Value = 723.17
$ws.Range("AB3").Value = 105.98
$ws.Range("AC3").Value = 498
$ws.Range("AD3").Value = 18.59
$ws.Range("AE3").Value = 10229
$ws.Range("AF3").Value = 0.91
$ws.Range("AG3").Value = 330
$ws.Range("AH3").Value = 3.56
$ws.Range("AI3").Value = 65.40000000000001
$ws.Range("AJ3").Value = 326700891

# --- Row 4 ---
$ws.Range("D4").Value = 78340
$ws.Range("E4").Value = 32
$ws.Range("F4").Value = 32
$ws.Range("G4").Value = 209
$ws.Range("H4").Value = 157
$ws.Range("I4").Value = 393
$ws.Range("J4").Value = -236
$ws.Range("K4").Value = 572678
$ws.Range("L4").Value = 506020
$ws.Range("M4").Value = 66658
$ws.Range("N4").Value = 66389
$ws.Range("O4").Value = 269
$ws.Range("P4").Value = 34020
$ws.Range("Q4").Value = 50146
$ws.Range("R4").Value = -540
$ws.Range("S4").Value = -33635
$ws.Range("T4").Value = 526
$ws.Range("U4").ClearContents()
$ws.Range("V4").Value = 67478
$ws.Range("W4").Value = 0.04
$ws.Range("X4").Value = 0.2
$ws.Range("Y4").Value = 0.78
$ws.Range("Z4").Value = 0.09
$ws.Range("AA4").Value = 759.13
$ws.Range("AB4").Value = 98.56999999999999
$ws.Range("AC4").Value = 115
$ws.Range("AD4").Value = 61.72
$ws.Range("AE4").Value = 12821
$ws.Range("AF4").Value = 0.55
$ws.Range("AG4").Value = 50
$ws.Range("AH4").Value = 0.71
$ws.Range("AI4").Value = 66.02
$ws.Range("AJ4").Value = 666316408

# --- Row 5 ---
$ws.Range("D5").Value = 102986
$ws.Range("E5").Value = 6278
$ws.Range("F5").Value = 6278
$ws.Range("G5").Value = 6647
$ws.Range("H5").Value = 5049
$ws.Range("I5").Value = 5032
$ws.Range("J5").Value = 18
$ws.Range("K5").Value = 906266
$ws.Range("L5").Value = 832421
$ws.Range("M5").Value = 73845
$ws.Range("N5").Value = 73567
$ws.Range("O5").Value = 278
$ws.Range("P5").Value = 34020
$ws.Range("Q5").Value = -41446
$ws.Range("R5").Value = -4586
$ws.Range("S5").Value = 40500
$ws.Range("T5").Value = 466
$ws.Range("U5").ClearContents()
$ws.Range("V5").Value = 68976
$ws.Range("W5").Value = 6.09
$ws.Range("X5").Value = 4.9
$ws.Range("Y5").Value = 7.19
$ws.Range("Z5").Value = 0.68
$ws.Range("AA5").Value = 1127.25
$ws.Range("AB5").Value = 118.9
$ws.Range("AC5").Value = 740
$ws.Range("AD5").Value = 12.15
$ws.Range("AE5").Value = 13002
$ws.Range("AF5").Value = 0.6899999999999999
$ws.Range("AG5").Value = 220
$ws.Range("AH5").Value = 2.45
$ws.Range("AI5").Value = 24.78
$ws.Range("AJ5").Value = 666316408

# --- Row 6 ---
$ws.Range("D6").Value = 133219
$ws.Range("E6").Value = 5123
$ws.Range("F6").Value = 5123
$ws.Range("G6").Value = 5850
$ws.Range("H6").Value = 4620
$ws.Range("I6").Value = 4566
$ws.Range("K6").Value = 1180855
$ws.Range("L6").Value = 1097331
$ws.Range("M6").Value = 83524
$ws.Range("N6").Value = 83136
$ws.Range("P6").Value = 41020
$ws.Range("Q6").Value = -48755
$ws.Range("R6").Value = -28419
$ws.Range("S6").Value = 80542
$ws.Range("T6").Value = 642
$ws.Range("U6").ClearContents()
$ws.Range("V6").Value = 102652
$ws.Range("W6").Value = 3.85
$ws.Range("X6").Value = 3.47
$ws.Range("Y6").Value = 5.83
$ws.Range("Z6").Value = 0.44
$ws.Range("AA6").Value = 1313.8
$ws.Range("AB6").Value = 105.13
$ws.Range("AC6").Value = 576
$ws.Range("AD6").Value = 11.34
$ws.Range("AE6").Value = 11900
$ws.Range("AF6").Value = 0.55
$ws.Range("AG6").Value = 220
$ws.Range("AH6").Value = 3.37
$ws.Range("AI6").Value = 33.71
$ws.Range("AJ6").Value = 658316408

# --- Row 7 ---
$ws.Range("D7").Value = 163672
$ws.Range("E7").Value = 7286
$ws.Range("G7").Value = 8676
$ws.Range("H7").Value = 6403
$ws.Range("I7").Value = 6368
$ws.Range("K7").Value = 1359433
$ws.Range("L7").Value = 1268468
$ws.Range("M7").Value = 90964
$ws.Range("N7").Value = 89884
$ws.Range("P7").Value = 41020
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").Value = 4.45
$ws.Range("X7").Value = 3.91
$ws.Range("Y7").Value = 7.36
$ws.Range("Z7").Value = 0.5
$ws.Range("AA7").Value = 1394.47
$ws.Range("AC7").Value = 784
$ws.Range("AD7").Value = 9.640000000000001
$ws.Range("AE7").Value = 12853
$ws.Range("AF7").Value = 0.59
$ws.Range("AG7").Value = 255
$ws.Range("AH7").Value = 3.38
$ws.Range("AI7").Value = 26.38

# --- Row 8 ---
$ws.Range("D8").Value = 147931
$ws.Range("E8").Value = 7598
$ws.Range("G8").Value = 8434
$ws.Range("H8").Value = 6235
$ws.Range("I8").Value = 6033
$ws.Range("K8").Value = 1441383
$ws.Range("L8").Value = 1345912
$ws.Range("M8").Value = 95469
$ws.Range("N8").Value = 94379
$ws.Range("P8").Value = 41020
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").Value = 5.14
$ws.Range("X8").Value = 4.21
$ws.Range("Y8").Value = 6.55
$ws.Range("Z8").Value = 0.45
$ws.Range("AA8").Value = 1409.79
$ws.Range("AC8").Value = 743
$ws.Range("AD8").Value = 9.279999999999999
$ws.Range("AE8").Value = 13496
$ws.Range("AF8").Value = 0.51
$ws.Range("AG8").Value = 265
$ws.Range("AH8").Value = 3.84
$ws.Range("AI8").Value = 28.9

# --- Row 9 ---
$ws.Range("D9").Value = 172710
$ws.Range("E9").Value = 8254
$ws.Range("G9").Value = 8772
$ws.Range("H9").Value = 6439
$ws.Range("I9").Value = 6349
$ws.Range("K9").Value = 1516386
$ws.Range("L9").Value = 1416326
$ws.Range("M9").Value = 99874
$ws.Range("N9").Value = 99332
$ws.Range("P9").Value = 41020
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").Value = 4.78
$ws.Range("X9").Value = 3.73
$ws.Range("Y9").Value = 6.55
$ws.Range("Z9").Value = 0.43
$ws.Range("AA9").Value = 1418.12
$ws.Range("AC9").Value = 782
$ws.Range("AD9").Value = 8.82
$ws.Range("AE9").Value = 14204
$ws.Range("AF9").Value = 0.49
$ws.Range("AG9").Value = 271
$ws.Range("AH9").Value = 3.94
$ws.Range("AI9").Value = 28.14
